$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.448.13'
$ws.Range("E2").Value = '  -0.94%  '
$ws.Range("D3").Value = '2.051.63'
$ws.Range("E3").Value = '  -1.77%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.28%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '228.97'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.10%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.614'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.74%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '56.37'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.22%  '
$ws.Range("E9").Value = '  -1.74%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0803'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.75%  '
$ws.Range("E11").Value = '  -2.01%  '
$ws.Range("D12").Value = '2.355.13'
$ws.Range("E12").Value = '  -1.70%  '
$ws.Range("E13").Value = '  -4.74%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.68'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.44%  '
$ws.Range("E15").Value = '  -3.03%  '
$ws.Range("E16").Value = '  -2.16%  '
$ws.Range("D17").Value = '2.048.12'
$ws.Range("E17").Value = '  -1.95%  '
$ws.Range("D18").Value = '37.338.08'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.06'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.03%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '69.78'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.69%  '
$ws.Range("D21").Value = '0.0₃0849'
$ws.Range("E21").Value = '  +1.71%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '225.75'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.76%  '
$ws.Range("E23").Value = '  +0.06%  '
$ws.Range("E24").Value = '  -1.18%  '
$ws.Range("E25").Value = '  -4.59%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.52'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.11%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '168.08'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.05%  '
$ws.Range("E28").Value = '  -5.18%  '
$ws.Range("E29").Value = '  -2.11%  '
$ws.Range("E30").Value = '  -2.85%  '
$ws.Range("E31").Value = '  -2.18%  '
$ws.Range("E32").Value = '  -3.34%  '
$ws.Range("E33").Value = '  -3.19%  '
$ws.Range("E34").Value = '  -2.24%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.39'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.13%  '
$ws.Range("E36").Value = '  +0.27%  '
$ws.Range("E37").Value = '  -0.09%  '
$ws.Range("E38").Value = '  -3.84%  '
$ws.Range("E39").Value = '  +0.81%  '
$ws.Range("E40").Value = '  -7.00%  '
$ws.Range("D41").Value = '1.494.02'
$ws.Range("E41").Value = '  +3.12%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '16.77'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.11%  '
$ws.Range("B43").Value = 'Cronos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0939'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.14%  '
$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '96.04'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -5.20%  '
$ws.Range("B45").Value = 'HuobiToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.86'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.00%  '
$ws.Range("E46").Value = '  -3.70%  '
$ws.Range("E47").Value = '  -4.18%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.21'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.13%  '
$ws.Range("E49").Value = '  -1.14%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.81'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -7.47%  '
$ws.Range("D51").Value = '2.240.08'
$ws.Range("E51").Value = '  -1.72%  '
